$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A4 values
$ws.Range("A2").Value = 24
$ws.Range("A3").Value = 37
$ws.Range("A4").Value = 40

# Update the selection to A5
$ws.Range("A5").Select()
